$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 19 new blank rows starting at row 4 -------------------------
# Inserting at row 4 repeatedly pushes the existing rows 4-10 further down
# and, crucially, each freshly inserted row inherits the formatting of the
# row immediately above it (so column A gets style "2" - the same style
# already used by rows 2-3 - without creating a brand-new style record).
for ($i = 0; $i -lt 19; $i++) {
    $ws.Rows.Item(4).Insert()
}

# --- Apartment codes (column A) -----------------------------------------
$codes = @(
    "A001", "A002", "A003", "A004", "A005", "A006", "A007", "A008", "A009",
    "A010", "A011", "A012", "A013", "A014", "A015", "A016", "A017", "A018",
    "A019", "A020"
)

# --- Usage (column B) ----------------------------------------------------
$usage = @(500, 344, 123, 542, 123, 455, 234, 543, 121,
           500, 344, 123, 542, 123, 455, 234, 543, 121, 543, 121)

for ($i = 0; $i -lt 20; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $codes[$i]
    $ws.Cells.Item($row, 2).Value = $usage[$i]
    $ws.Cells.Item($row, 3).Value = "kWh"
}

# --- Row 22: leftover blank row, column A only ---------------------------
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()

# --- Drop the old (now pushed-down and unused) source rows ---------------
$ws.Rows.Item("23:29").Delete()

# --- Selection, matching the author's final cursor position --------------
$ws.Range("G15").Select()
